$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 43 from
# 2025-06-22 (serial 45830) to 2025-06-23 (serial 45831).
for ($r = 2; $r -le 43; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45830) {
        $cell.Value2 = 45831
    }
}
